{"js": "// Update the multiplication-table answer sheet: replace each old\n// \"A\u00d7B=C\" equation text with its new value (25 distinct substitutions).\nconst replacements = [\n  [\"60\u00d749=2940\", \"52\u00d721=1092\"],\n  [\"80\u00d795=7600\", \"68\u00d796=6528\"],\n  [\"49\u00d755=2695\", \"31\u00d771=2201\"],\n  [\"64\u00d732=2048\", \"76\u00d740=3040\"],\n  [\"99\u00d712=1188\", \"83\u00d762=5146\"],\n  [\"50\u00d782=4100\", \"67\u00d788=5896\"],\n  [\"69\u00d796=6624\", \"61\u00d767=4087\"],\n  [\"66\u00d718=1188\", \"40\u00d721=840\"],\n  [\"90\u00d783=7470\", \"17\u00d793=1581\"],\n  [\"82\u00d712=984\", \"15\u00d794=1410\"],\n  [\"41\u00d715=615\", \"72\u00d780=5760\"],\n  [\"24\u00d714=336\", \"57\u00d753=3021\"],\n  [\"74\u00d746=3404\", \"28\u00d738=1064\"],\n  [\"14\u00d780=1120\", \"75\u00d775=5625\"],\n  [\"71\u00d785=6035\", \"62\u00d717=1054\"],\n  [\"42\u00d767=2814\", \"63\u00d782=5166\"],\n  [\"34\u00d754=1836\", \"23\u00d765=1495\"],\n  [\"13\u00d719=247\", \"90\u00d793=8370\"],\n  [\"44\u00d756=2464\", \"84\u00d721=1764\"],\n  [\"80\u00d745=3600\", \"11\u00d717=187\"],\n  [\"70\u00d735=2450\", \"90\u00d720=1800\"],\n  [\"74\u00d737=2738\", \"93\u00d729=2697\"],\n  [\"22\u00d725=550\", \"75\u00d795=7125\"],\n  [\"53\u00d784=4452\", \"94\u00d735=3290\"],\n  [\"40\u00d761=2440\", \"48\u00d732=1536\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-table answer sheet: replace each old\n# \"A\u00d7B=C\" equation text with its new value (25 distinct substitutions).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"60\u00d749=2940\", \"52\u00d721=1092\"),\n  @(\"80\u00d795=7600\", \"68\u00d796=6528\"),\n  @(\"49\u00d755=2695\", \"31\u00d771=2201\"),\n  @(\"64\u00d732=2048\", \"76\u00d740=3040\"),\n  @(\"99\u00d712=1188\", \"83\u00d762=5146\"),\n  @(\"50\u00d782=4100\", \"67\u00d788=5896\"),\n  @(\"69\u00d796=6624\", \"61\u00d767=4087\"),\n  @(\"66\u00d718=1188\", \"40\u00d721=840\"),\n  @(\"90\u00d783=7470\", \"17\u00d793=1581\"),\n  @(\"82\u00d712=984\", \"15\u00d794=1410\"),\n  @(\"41\u00d715=615\", \"72\u00d780=5760\"),\n  @(\"24\u00d714=336\", \"57\u00d753=3021\"),\n  @(\"74\u00d746=3404\", \"28\u00d738=1064\"),\n  @(\"14\u00d780=1120\", \"75\u00d775=5625\"),\n  @(\"71\u00d785=6035\", \"62\u00d717=1054\"),\n  @(\"42\u00d767=2814\", \"63\u00d782=5166\"),\n  @(\"34\u00d754=1836\", \"23\u00d765=1495\"),\n  @(\"13\u00d719=247\", \"90\u00d793=8370\"),\n  @(\"44\u00d756=2464\", \"84\u00d721=1764\"),\n  @(\"80\u00d745=3600\", \"11\u00d717=187\"),\n  @(\"70\u00d735=2450\", \"90\u00d720=1800\"),\n  @(\"74\u00d737=2738\", \"93\u00d729=2697\"),\n  @(\"22\u00d725=550\", \"75\u00d795=7125\"),\n  @(\"53\u00d784=4452\", \"94\u00d735=3290\"),\n  @(\"40\u00d761=2440\", \"48\u00d732=1536\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Could not find text to replace: $oldText\"\n  }\n}\n"}
